$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row-structure changes first, so row numbers below are final -----

# Insert a new row at row 11 ("app_launch"); everything from the old
# row 11 ("citrix_edit") downward shifts down by one.
$ws.Rows.Item(11).Insert()

# Insert three new rows before the row that now holds "view_create"
# (row 22 originally, row 23 after the previous insert) for the
# duplicated sound-task entries.
$ws.Rows.Item(23).Resize(3).Insert()

# --- Now fill in the new / changed cell values -------------------------

# The two "background" rows (now rows 43-44) become the new settings
# entries.
$ws.Range("A43").Value = "settings_background"
$ws.Range("A44").Value = "settings_theme"

# "import_export" (now row 47) gets an "N/A" result in column B.
$ws.Range("B47").Value = "N/A"

# Append five new wifi-related rows at the bottom of the table.
$ws.Range("A53").Value = "wifi_hidden_ssid"
$ws.Range("B53").Value = "N/A"
$ws.Range("A54").Value = "wifi_modify_icon"
$ws.Range("B54").Value = "N/A"
$ws.Range("A55").Value = "wifi_modify_settings"
$ws.Range("B55").Value = "N/A"
$ws.Range("A56").Value = "wifi_readonly"
$ws.Range("B56").Value = "N/A"
$ws.Range("A57").Value = "wifi_wpap"
$ws.Range("B57").Value = "N/A"

# Fill the three rows inserted earlier for the sound-task entries.
$ws.Range("A23").Value = "task_modify_sound_key"
$ws.Range("A24").Value = "task_modify_sound_mouse"
$ws.Range("A25").Value = "task_sound_readonly"

# Finally, the new row 11.
$ws.Range("A11").Value = "app_launch"

# Match the saved selection (the author's last-focused cell).
$ws.Range("B57").Select() | Out-Null
